$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Info": update Objetivo / Tiempo result values
# ---------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("A2").Value = 3183447082398.467
$wsInfo.Range("B2").Value = 2.180000066757202

# ---------------------------------------------------------------
# Sheet "Activados": extend data from 3 rows (A1:B4) to 19 rows
# (A1:B20), with Proceso = 1 and Tiempo stepping by 20
# ---------------------------------------------------------------
$wsAct = $wb.Worksheets.Item("Activados")
for ($i = 0; $i -lt 19; $i++) {
    $row = 2 + $i
    $wsAct.Cells.Item($row, 1).Value = 1
    $wsAct.Cells.Item($row, 2).Value = $i * 20
}

# ---------------------------------------------------------------
# Sheet "Operando": Proceso column changes from 4 to 1 for all
# 365 data rows (Tiempo column is unchanged)
# ---------------------------------------------------------------
$wsOp = $wb.Worksheets.Item("Operando")
$wsOp.Range("A2:A366").Value = 1

# ---------------------------------------------------------------
# Sheet "Contaminantes": update Z / Concentracion result values
# ---------------------------------------------------------------
$wsCont = $wb.Worksheets.Item("Contaminantes")
$wsCont.Cells.Item(2, 2).Value = 2650328644320
$wsCont.Cells.Item(2, 3).Value = 98.294
$wsCont.Cells.Item(3, 2).Value = 146949876000
$wsCont.Cells.Item(3, 3).Value = 5.449999999999999
$wsCont.Cells.Item(4, 2).Value = 124219830960
$wsCont.Cells.Item(4, 3).Value = 4.606999999999998
$wsCont.Cells.Item(5, 2).Value = 465918.4679472
$wsCont.Cells.Item(5, 3).Value = 0.00001727974
$wsCont.Cells.Item(6, 2).Value = 261948265200
$wsCont.Cells.Item(6, 3).Value = 9.715000000000002
